$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 26: update title (D26)
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 46: update title (D46) and link (E46)
$ws.Range("D46").Value = "BIRADS (유방영상보고 및 자료체계)"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/438"

# Row 51: update title (D51) and link (E51)
$ws.Range("D51").Value = "[flask+jinja2] 반올림하기, round 필터"
$ws.Range("E51").Value = "https://bskyvision.com/1228"
